$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1424555306933115
$ws.Cells.Item(2, 4).Value = 0.01226424981397045
$ws.Cells.Item(2, 5).Value = 0.4190205312761037
$ws.Cells.Item(2, 6).Value = 0.8117438529910004
$ws.Cells.Item(2, 7).Value = 0.669986943447725
$ws.Cells.Item(2, 8).Value = 0.7100009786584565
$ws.Cells.Item(2, 11).Value = 0.932193313258125
$ws.Cells.Item(2, 14).Value = 1.293435771131678
$ws.Cells.Item(3, 2).Value = 0.1329270018370465
$ws.Cells.Item(3, 4).Value = 0.01222692433373318
$ws.Cells.Item(3, 5).Value = 0.3653875206293691
$ws.Cells.Item(3, 6).Value = 0.7910567660141226
$ws.Cells.Item(3, 7).Value = 0.6484228125779907
$ws.Cells.Item(3, 8).Value = 0.7055526869364428
$ws.Cells.Item(3, 11).Value = 0.8130736497739974
$ws.Cells.Item(3, 14).Value = 1.311679579673037
$ws.Cells.Item(4, 2).Value = 0.1271470273094053
$ws.Cells.Item(4, 4).Value = 0.01220573269201353
$ws.Cells.Item(4, 5).Value = 0.3325690400924088
$ws.Cells.Item(4, 6).Value = 0.7790000758484581
$ws.Cells.Item(4, 7).Value = 0.6357724868773289
$ws.Cells.Item(4, 8).Value = 0.7032943838607366
$ws.Cells.Item(4, 11).Value = 0.7399104512414283
$ws.Cells.Item(4, 14).Value = 1.323435997925484
$ws.Cells.Item(5, 2).Value = 0.124809526523407
$ws.Cells.Item(5, 4).Value = 0.01219753319307415
$ws.Cells.Item(5, 5).Value = 0.3192211660014124
$ws.Cells.Item(5, 6).Value = 0.7742481333345381
$ws.Cells.Item(5, 7).Value = 0.6307645355264242
$ws.Cells.Item(5, 8).Value = 0.7024925944058396
$ws.Cells.Item(5, 11).Value = 0.7100891770322448
$ws.Cells.Item(5, 14).Value = 1.328366122324325
$ws.Cells.Item(6, 2).Value = 0.1244224706180006
$ws.Cells.Item(6, 4).Value = 0.01219619809329053
$ws.Cells.Item(6, 5).Value = 0.3170062593404595
$ws.Cells.Item(6, 6).Value = 0.7734687873153661
$ws.Cells.Item(6, 7).Value = 0.629941821098015
$ws.Cells.Item(6, 8).Value = 0.7023666013440533
$ws.Cells.Item(6, 11).Value = 0.7051369501528484
$ws.Cells.Item(6, 14).Value = 1.329193173161372
$ws.Cells.Item(7, 2).Value = 0.1271154303542943
$ws.Cells.Item(7, 4).Value = 0.01220562034221295
$ws.Cells.Item(7, 5).Value = 0.3323889240920721
$ws.Cells.Item(7, 6).Value = 0.7789353378959873
$ws.Cells.Item(7, 7).Value = 0.6357043536895901
$ws.Cells.Item(7, 8).Value = 0.7032830914920822
$ws.Cells.Item(7, 11).Value = 0.7395082992187838
$ws.Cells.Item(7, 14).Value = 1.323501923583007
$ws.Cells.Item(8, 2).Value = 0.139155524127716
$ws.Cells.Item(8, 4).Value = 0.01225102281800083
$ws.Cells.Item(8, 5).Value = 0.4005031094127531
$ws.Cells.Item(8, 6).Value = 0.8044764077536541
$ws.Cells.Item(8, 7).Value = 0.6624283935808819
$ws.Cells.Item(8, 8).Value = 0.7083687717010179
$ws.Cells.Item(8, 11).Value = 0.8911251354244882
$ws.Cells.Item(8, 14).Value = 1.299611011482763
$ws.Cells.Item(9, 2).Value = 0.1633214294895282
$ws.Cells.Item(9, 4).Value = 0.01235369287405064
$ws.Cells.Item(9, 5).Value = 0.5350789449597642
$ws.Cells.Item(9, 6).Value = 0.8597313139756864
$ws.Cells.Item(9, 7).Value = 0.719576149890429
$ws.Cells.Item(9, 8).Value = 0.722116260616815
$ws.Cells.Item(9, 11).Value = 1.188313826996136
$ws.Cells.Item(9, 14).Value = 1.257169298009467
$ws.Cells.Item(10, 2).Value = 0.1814103069771704
$ws.Cells.Item(10, 4).Value = 0.01243740363141299
$ws.Cells.Item(10, 5).Value = 0.6347289413529609
$ws.Cells.Item(10, 6).Value = 0.90355182278077
$ws.Cells.Item(10, 7).Value = 0.7645413002500732
$ws.Cells.Item(10, 8).Value = 0.7345500852917439
$ws.Cells.Item(10, 11).Value = 1.406673835499873
$ws.Cells.Item(10, 14).Value = 1.228683954465931
$ws.Cells.Item(11, 2).Value = 0.1897111742345459
$ws.Cells.Item(11, 4).Value = 0.01247728789495639
$ws.Cells.Item(11, 5).Value = 0.6802676102145426
$ws.Cells.Item(11, 6).Value = 0.9242034884202184
$ws.Cells.Item(11, 7).Value = 0.7856633329674878
$ws.Cells.Item(11, 8).Value = 0.7407205230536533
$ws.Cells.Item(11, 11).Value = 1.506037106683891
$ws.Cells.Item(11, 14).Value = 1.216312791333018
$ws.Cells.Item(12, 2).Value = 0.1928647603465521
$ws.Cells.Item(12, 4).Value = 0.01249265121182574
$ws.Cells.Item(12, 5).Value = 0.6975445672242415
$ws.Cells.Item(12, 6).Value = 0.9321281463112143
$ws.Cells.Item(12, 7).Value = 0.793759122868579
$ws.Cells.Item(12, 8).Value = 0.7431316107924602
$ws.Cells.Item(12, 11).Value = 1.543669102161118
$ws.Cells.Item(12, 14).Value = 1.211712773817323
$ws.Cells.Item(13, 2).Value = 0.1921851265993411
$ws.Cells.Item(13, 4).Value = 0.0124893308663836
$ws.Cells.Item(13, 5).Value = 0.6938221798461086
$ws.Cells.Item(13, 6).Value = 0.93041677139
$ws.Cells.Item(13, 7).Value = 0.7920111982661808
$ws.Cells.Item(13, 8).Value = 0.7426090189811703
$ws.Cells.Item(13, 11).Value = 1.535564128501392
$ws.Cells.Item(13, 14).Value = 1.212699699669943
$ws.Cells.Item(14, 2).Value = 0.1899704175196177
$ws.Cells.Item(14, 4).Value = 0.01247854662842229
$ws.Cells.Item(14, 5).Value = 0.6816883302889778
$ws.Cells.Item(14, 6).Value = 0.9248533574533155
$ws.Cells.Item(14, 7).Value = 0.7863274196880923
$ws.Cells.Item(14, 8).Value = 0.7409173889944043
$ws.Cells.Item(14, 11).Value = 1.509133002286774
$ws.Cells.Item(14, 14).Value = 1.21593264658021
$ws.Cells.Item(15, 2).Value = 0.1886151718771174
$ws.Cells.Item(15, 4).Value = 0.0124719748593769
$ws.Cells.Item(15, 5).Value = 0.6742603019180109
$ws.Cells.Item(15, 6).Value = 0.9214592255745373
$ws.Cells.Item(15, 7).Value = 0.7828586587539235
$ws.Cells.Item(15, 8).Value = 0.7398909317726634
$ws.Cells.Item(15, 11).Value = 1.492943890367201
$ws.Cells.Item(15, 14).Value = 1.217923954979147
$ws.Cells.Item(16, 2).Value = 0.1808692649986625
$ws.Cells.Item(16, 4).Value = 0.01243483344053686
$ws.Cells.Item(16, 5).Value = 0.6317572792954138
$ws.Cells.Item(16, 6).Value = 0.9022167233194125
$ws.Cells.Item(16, 7).Value = 0.7631744698508385
$ws.Cells.Item(16, 8).Value = 0.7341572255183735
$ws.Cells.Item(16, 11).Value = 1.400180899380018
$ws.Cells.Item(16, 14).Value = 1.229504273115646
$ws.Cells.Item(17, 2).Value = 0.1761357759129822
$ws.Cells.Item(17, 4).Value = 0.01241251065962246
$ws.Cells.Item(17, 5).Value = 0.605738027076967
$ws.Cells.Item(17, 6).Value = 0.8905966501634168
$ws.Cells.Item(17, 7).Value = 0.7512707231278171
$ws.Cells.Item(17, 8).Value = 0.7307718644117358
$ws.Cells.Item(17, 11).Value = 1.343282203908643
$ws.Cells.Item(17, 14).Value = 1.236758947891944
$ws.Cells.Item(18, 2).Value = 0.1734199976769304
$ws.Cells.Item(18, 4).Value = 0.01239984090144475
$ws.Cells.Item(18, 5).Value = 0.5907918019628084
$ws.Cells.Item(18, 6).Value = 0.8839805563274723
$ws.Cells.Item(18, 7).Value = 0.7444867146987519
$ws.Cells.Item(18, 8).Value = 0.7288730636407763
$ws.Cells.Item(18, 11).Value = 1.310558364708186
$ws.Cells.Item(18, 14).Value = 1.24098687620472
$ws.Cells.Item(19, 2).Value = 0.1725016541232804
$ws.Cells.Item(19, 4).Value = 0.01239558027650034
$ws.Cells.Item(19, 5).Value = 0.5857345126585614
$ws.Cells.Item(19, 6).Value = 0.8817520131094341
$ws.Cells.Item(19, 7).Value = 0.7422004944750142
$ws.Cells.Item(19, 8).Value = 0.7282384537403175
$ws.Cells.Item(19, 11).Value = 1.299479085360758
$ws.Cells.Item(19, 14).Value = 1.24242785799822
$ws.Cells.Item(20, 2).Value = 0.1766389612490542
$ws.Cells.Item(20, 4).Value = 0.01241486939111169
$ws.Cells.Item(20, 5).Value = 0.6085057981708388
$ws.Cells.Item(20, 6).Value = 0.8918266361928602
$ws.Cells.Item(20, 7).Value = 0.752531396787532
$ws.Cells.Item(20, 8).Value = 0.7311272320291948
$ws.Cells.Item(20, 11).Value = 1.349338872160729
$ws.Cells.Item(20, 14).Value = 1.235980957248712
$ws.Cells.Item(21, 2).Value = 0.1906206548486438
$ws.Cells.Item(21, 4).Value = 0.01248170715631147
$ws.Cells.Item(21, 5).Value = 0.6852514341639591
$ws.Cells.Item(21, 6).Value = 0.9264846265690778
$ws.Cells.Item(21, 7).Value = 0.7879942303894438
$ws.Cells.Item(21, 8).Value = 0.7414122361719819
$ws.Cells.Item(21, 11).Value = 1.51689631381538
$ws.Cells.Item(21, 14).Value = 1.214980751517651
$ws.Cells.Item(22, 2).Value = 0.1998180476458202
$ws.Cells.Item(22, 4).Value = 0.01252690570665749
$ws.Cells.Item(22, 5).Value = 0.7355995436218024
$ws.Cells.Item(22, 6).Value = 0.9497442264093081
$ws.Cells.Item(22, 7).Value = 0.8117392415180404
$ws.Cells.Item(22, 8).Value = 0.748568441317218
$ws.Cells.Item(22, 11).Value = 1.626436026247347
$ws.Cells.Item(22, 14).Value = 1.201749573893719
$ws.Cells.Item(23, 2).Value = 0.1949038257850191
$ws.Cells.Item(23, 4).Value = 0.01250264331291362
$ws.Cells.Item(23, 5).Value = 0.7087095095375844
$ws.Cells.Item(23, 6).Value = 0.9372740714552208
$ws.Cells.Item(23, 7).Value = 0.7990136445334883
$ws.Cells.Item(23, 8).Value = 0.7447091183319401
$ws.Cells.Item(23, 11).Value = 1.567969444713867
$ws.Cells.Item(23, 14).Value = 1.208766041080379
$ws.Cells.Item(24, 2).Value = 0.1764114538052723
$ws.Cells.Item(24, 4).Value = 0.01241380249772561
$ws.Cells.Item(24, 5).Value = 0.6072544498706378
$ws.Cells.Item(24, 6).Value = 0.8912703588659667
$ws.Cells.Item(24, 7).Value = 0.7519612606392343
$ws.Cells.Item(24, 8).Value = 0.7309664224420658
$ws.Cells.Item(24, 11).Value = 1.346600689970387
$ws.Cells.Item(24, 14).Value = 1.236332508912076
$ws.Cells.Item(25, 2).Value = 0.156724912090823
$ws.Cells.Item(25, 4).Value = 0.01232446347498417
$ws.Cells.Item(25, 5).Value = 0.4985481870050847
$ws.Cells.Item(25, 6).Value = 0.8442224631986761
$ws.Cells.Item(25, 7).Value = 0.7035994734058875
$ws.Cells.Item(25, 8).Value = 0.7179897809866986
$ws.Cells.Item(25, 11).Value = 1.107920699135832
$ws.Cells.Item(25, 14).Value = 1.268178042075494
